$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1580
$ws1.Range("F3").Value = 3322
$ws1.Range("F4").Value = 28
$ws1.Range("F5").Value = 753
$ws1.Range("F6").Value = 2349
$ws1.Range("F7").Value = 507
$ws1.Range("F8").Value = 424
$ws1.Range("F9").Value = 254
$ws1.Range("F11").Value = 370
$ws1.Range("F12").Value = 1112
$ws1.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202406/ckVVTuNj1717752114555.jpeg"
$ws1.Range("F15").Value = 92
$ws1.Range("F16").Value = 273
$ws1.Range("F17").Value = 4875
$ws1.Range("F18").Value = 29
$ws1.Range("F19").Value = 1378
$ws1.Range("F20").Value = 3578
$ws1.Range("F21").Value = 144
$ws1.Range("F22").Value = 206
$ws1.Range("F23").Value = 3876
$ws1.Range("F24").Value = 5244
$ws1.Range("F27").Value = 575
$ws1.Range("F28").Value = 3373
$ws1.Range("F29").Value = 390
$ws1.Range("F31").Value = 149
$ws1.Range("F34").Value = 1220
$ws1.Range("F35").Value = 33
$ws1.Range("F36").Value = 51
$ws1.Range("F37").Value = 1444
$ws1.Range("F39").Value = 1426
$ws1.Range("F40").Value = 35
$ws1.Range("F41").Value = 921
$ws1.Range("F42").Value = 899
$ws1.Range("F43").Value = 527
$ws1.Range("F45").Value = 2114
$ws1.Range("F46").Value = 88
$ws1.Range("F47").Value = 185
$ws1.Range("F48").Value = 374
$ws1.Range("F49").Value = 3760

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 1027
$ws2.Range("G6").Value = 319

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 2481

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 2481
$ws4.Range("F3").Value = 1580
$ws4.Range("F4").Value = 3322
$ws4.Range("F5").Value = 28
$ws4.Range("F6").Value = 753
$ws4.Range("F8").Value = 2349
$ws4.Range("F9").Value = 507
$ws4.Range("F10").Value = 424
$ws4.Range("F11").Value = 254
$ws4.Range("F12").Value = 1027
$ws4.Range("G12").Value = 319
$ws4.Range("F14").Value = 370
$ws4.Range("F15").Value = 1112
$ws4.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202406/ckVVTuNj1717752114555.jpeg"
$ws4.Range("F18").Value = 92
$ws4.Range("F19").Value = 273
$ws4.Range("F20").Value = 4875
$ws4.Range("F22").Value = 1378
$ws4.Range("F23").Value = 3876
$ws4.Range("F24").Value = 5244
$ws4.Range("F27").Value = 575
$ws4.Range("F28").Value = 3373
$ws4.Range("F29").Value = 390
$ws4.Range("F31").Value = 149
$ws4.Range("F33").Value = 1220
$ws4.Range("F34").Value = 33
$ws4.Range("F35").Value = 51
$ws4.Range("F36").Value = 1444
$ws4.Range("F37").Value = 1426
$ws4.Range("F38").Value = 921
$ws4.Range("F39").Value = 527
$ws4.Range("F43").Value = 2118
$ws4.Range("F45").Value = 88
$ws4.Range("F46").Value = 185
$ws4.Range("F47").Value = 374
$ws4.Range("F49").Value = 3760
